$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new quarter columns R, S, T (31/12/2023, 31/03/2024, 30/06/2024)
# Match style of existing header cells (bold, centered, bordered) via direct formatting
$ws.Range("R1:T1").Font.Bold = $true
$ws.Range("R1:T1").HorizontalAlignment = -4108
$ws.Range("R1:T1").VerticalAlignment = -4160
$ws.Range("R1:T1").Borders.LineStyle = 1
$ws.Range("R1").Value = "31/12/2023"
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

$ws.Range("R2").Value = 2421372.928
$ws.Range("S2").Value = 2451693.056
$ws.Range("T2").Value = 2455854.08
$ws.Range("R3").Value = 1991437.056
$ws.Range("S3").Value = 1938444.032
$ws.Range("T3").Value = 2057707.008
$ws.Range("R4").Value = 1580
$ws.Range("S4").Value = 1045
$ws.Range("T4").Value = 812
$ws.Range("R5").Value = 397595.008
$ws.Range("S5").Value = 447863.008
$ws.Range("T5").Value = 499878.016
$ws.Range("R6").Value = 1592262.016
$ws.Range("S6").Value = 1489536
$ws.Range("T6").Value = 1557016.96
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("R12").Value = 395336
$ws.Range("S12").Value = 476065.984
$ws.Range("T12").Value = 360188.992
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("R14").Value = 1039
$ws.Range("S14").Value = 5981
$ws.Range("T14").Value = 5659
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("R16").Value = 394075.008
$ws.Range("S16").Value = 469863.008
$ws.Range("T16").Value = 353720.992
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("T19").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("R21").Value = 222
$ws.Range("S21").Value = 222
$ws.Range("T21").Value = 809
$ws.Range("R22").Value = 23385
$ws.Range("S22").Value = 23627
$ws.Range("T22").Value = 24941
$ws.Range("R23").Value = 11086
$ws.Range("S23").Value = 13457
$ws.Range("T23").Value = 12808
$ws.Range("R24").Value = 129
$ws.Range("S24").Value = 99
$ws.Range("T24").Value = 209
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = 0
$ws.Range("T25").Value = 0
$ws.Range("R26").Value = 2421372.928
$ws.Range("S26").Value = 2451693.056
$ws.Range("T26").Value = 2455854.08
$ws.Range("R27").Value = 685825.024
$ws.Range("S27").Value = 660094.0159999999
$ws.Range("T27").Value = 721454.976
$ws.Range("R28").Value = 0
$ws.Range("S28").Value = 0
$ws.Range("T28").Value = 0
$ws.Range("R29").Value = 53465
$ws.Range("S29").Value = 41188
$ws.Range("T29").Value = 33947
$ws.Range("R30").Value = 27953
$ws.Range("S30").Value = 26588
$ws.Range("T30").Value = 21378
$ws.Range("R31").Value = 145523.008
$ws.Range("S31").Value = 115564
$ws.Range("T31").Value = 160138
$ws.Range("R32").Value = 43406
$ws.Range("S32").Value = 9944
$ws.Range("T32").Value = 6838
$ws.Range("R33").Value = 0
$ws.Range("S33").Value = 0
$ws.Range("T33").Value = 0
$ws.Range("R34").Value = 356720.992
$ws.Range("S34").Value = 407576
$ws.Range("T34").Value = 436316.992
$ws.Range("R35").Value = 58757
$ws.Range("S35").Value = 59234
$ws.Range("T35").Value = 62837
$ws.Range("R36").Value = 0
$ws.Range("S36").Value = 0
$ws.Range("T36").Value = 0
$ws.Range("R37").Value = 410844
$ws.Range("S37").Value = 423398.016
$ws.Range("T37").Value = 447100
$ws.Range("R38").Value = 139736
$ws.Range("S38").Value = 172914
$ws.Range("T38").Value = 160796
$ws.Range("R39").Value = 0
$ws.Range("S39").Value = 0
$ws.Range("T39").Value = 0
$ws.Range("R40").Value = 220508
$ws.Range("S40").Value = 202778
$ws.Range("T40").Value = 236743.008
$ws.Range("R41").Value = 16807
$ws.Range("S41").Value = 17409
$ws.Range("T41").Value = 21799
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 0
$ws.Range("T42").Value = 0
$ws.Range("R43").Value = 33793
$ws.Range("S43").Value = 30297
$ws.Range("T43").Value = 27762
$ws.Range("R44").Value = 0
$ws.Range("S44").Value = 0
$ws.Range("T44").Value = 0
$ws.Range("R45").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("T45").Value = 0
$ws.Range("R46").Value = 119507
$ws.Range("S46").Value = 131545
$ws.Range("T46").Value = 118377
$ws.Range("R47").Value = 1205197
$ws.Range("S47").Value = 1236655.96
$ws.Range("T47").Value = 1168921.944
$ws.Range("R48").Value = 1109028.992
$ws.Range("S48").Value = 1109028.992
$ws.Range("T48").Value = 1109028.992
$ws.Range("R49").Value = -4641
$ws.Range("S49").Value = -3522
$ws.Range("T49").Value = -1065
$ws.Range("R50").Value = 0
$ws.Range("S50").Value = 0
$ws.Range("T50").Value = 0
$ws.Range("R51").Value = 100809
$ws.Range("S51").Value = 100809
$ws.Range("T51").Value = 59398
$ws.Range("R52").Value = 0
$ws.Range("S52").Value = 30340
$ws.Range("T52").Value = 1560
$ws.Range("R53").Value = 0
$ws.Range("S53").Value = 0
$ws.Range("T53").Value = 0
$ws.Range("R54").Value = 0
$ws.Range("S54").Value = 0
$ws.Range("T54").Value = 0
$ws.Range("R55").Value = 0
$ws.Range("S55").Value = 0
$ws.Range("T55").Value = 0
$ws.Range("R56").Value = 0
$ws.Range("S56").Value = 0
$ws.Range("T56").Value = 0
$ws.Range("R59").Value = 384420.032
$ws.Range("S59").Value = 245143.008
$ws.Range("T59").Value = 140815.008
$ws.Range("R60").Value = -280841.984
$ws.Range("S60").Value = -191020.992
$ws.Range("T60").Value = -135174
$ws.Range("R61").Value = 103578.008
$ws.Range("S61").Value = 54122
$ws.Range("T61").Value = 5641
$ws.Range("R62").Value = -18272
$ws.Range("S62").Value = -14634
$ws.Range("T62").Value = -15184
$ws.Range("R63").Value = -18749
$ws.Range("S63").Value = -15543
$ws.Range("T63").Value = -19294
$ws.Range("R64").Value = 4312
$ws.Range("S64").Value = 205
$ws.Range("T64").Value = 700
$ws.Range("R65").Value = -4183
$ws.Range("S65").Value = 457
$ws.Range("T65").Value = 0
$ws.Range("R66").Value = -5897
$ws.Range("S66").Value = 0
$ws.Range("T66").Value = -11932
$ws.Range("R67").Value = 3070
$ws.Range("S67").Value = -961
$ws.Range("T67").Value = 2425
$ws.Range("R68").Value = 23175
$ws.Range("S68").Value = 24214
$ws.Range("T68").Value = 20351
$ws.Range("R69").Value = 48864
$ws.Range("S69").Value = 28276
$ws.Range("T69").Value = 21204
$ws.Range("R70").Value = -25689
$ws.Range("S70").Value = -4062
$ws.Range("T70").Value = -853
$ws.Range("R74").Value = 87034
$ws.Range("S74").Value = 47860
$ws.Range("T74").Value = -17293
$ws.Range("R75").Value = -7871
$ws.Range("S75").Value = -4880
$ws.Range("T75").Value = -8130
$ws.Range("R76").Value = -2204
$ws.Range("S76").Value = -602
$ws.Range("T76").Value = 1559
$ws.Range("R79").Value = -35811
$ws.Range("S79").Value = -12038
$ws.Range("T79").Value = 0
$ws.Range("R80").Value = 41148
$ws.Range("S80").Value = 30340
$ws.Range("T80").Value = -23864

# Rows that are section headers with no data (blank placeholder cells), matching existing pattern
$ws.Range("R57:T58").Borders.LineStyle = 0
$ws.Range("R71:T73").Borders.LineStyle = 0
$ws.Range("R77:T78").Borders.LineStyle = 0
